$wb = $excel.ActiveWorkbook

# ---- Sheet ALC: 57 cell updates ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 25000
$ws.Range("J3").Value = 25000
$ws.Range("L3").Value = 25000
$ws.Range("N3").Value = -25228
$ws.Range("H11").Value = 75.72727
$ws.Range("I11").Value = 75.72727
$ws.Range("K11").Value = 75.72727
$ws.Range("M11").Value = 64.27273
$ws.Range("H19").Value = 1830.2307
$ws.Range("I19").Value = 1248.4
$ws.Range("J19").Value = 1968.762
$ws.Range("K19").Value = 1248.4
$ws.Range("L19").Value = 1968.762
$ws.Range("M19").Value = -1073.4
$ws.Range("N19").Value = -2318.762
$ws.Range("H70").Value = 2527.0908
$ws.Range("I70").Value = 4874.5
$ws.Range("J70").Value = 1185.7142
$ws.Range("K70").Value = 14623.5
$ws.Range("L70").Value = 3557.1426
$ws.Range("M70").Value = -14353.5
$ws.Range("N70").Value = -4097.142599999999
$ws.Range("H73").Value = 2527.0908
$ws.Range("I73").Value = 4874.5
$ws.Range("J73").Value = 1185.7142
$ws.Range("K73").Value = 14623.5
$ws.Range("L73").Value = 3557.1426
$ws.Range("M73").Value = -13687.5
$ws.Range("N73").Value = -5429.142599999999
$ws.Range("H76").Value = 12000
$ws.Range("I76").Value = 9000
$ws.Range("K76").Value = 9000
$ws.Range("M76").Value = -8685
$ws.Range("H79").Value = 12000
$ws.Range("I79").Value = 9000
$ws.Range("K79").Value = 9000
$ws.Range("M79").Value = -7908
$ws.Range("H95").Value = 16631.715
$ws.Range("J95").Value = 16631.715
$ws.Range("L95").Value = 16631.715
$ws.Range("N95").Value = -22123.715
$ws.Range("H102").Value = 25000
$ws.Range("J102").Value = 25000
$ws.Range("L102").Value = 25000
$ws.Range("N102").Value = -31490
$ws.Range("H137").Value = 1705
$ws.Range("I137").Value = 1441.5555
$ws.Range("K137").Value = 4324.666499999999
$ws.Range("M137").Value = -1774.666499999999
$ws.Range("H138").Value = 2051
$ws.Range("I138").Value = 1904.75
$ws.Range("K138").Value = 5714.25
$ws.Range("M138").Value = -574.25
$ws.Range("H141").Value = 3911.625
$ws.Range("I141").Value = 2327.5715
$ws.Range("K141").Value = 6982.7145
$ws.Range("M141").Value = -1802.7145

# ---- Sheet ARM: 27 cell updates ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 4058.6
$ws.Range("I13").Value = 1931.3334
$ws.Range("J13").Value = 7249.5
$ws.Range("K13").Value = 1931.3334
$ws.Range("L13").Value = 7249.5
$ws.Range("M13").Value = -1787.3334
$ws.Range("N13").Value = -7537.5
$ws.Range("H32").Value = 7994.9473
$ws.Range("I32").Value = 6349.4287
$ws.Range("K32").Value = 6349.4287
$ws.Range("M32").Value = -6062.4287
$ws.Range("H74").Value = 2976.4736
$ws.Range("I74").Value = 2603.625
$ws.Range("K74").Value = 2603.625
$ws.Range("M74").Value = -1729.625
$ws.Range("H77").Value = 2976.4736
$ws.Range("I77").Value = 2603.625
$ws.Range("K77").Value = 13018.125
$ws.Range("M77").Value = -8650.125
$ws.Range("H88").Value = 1389.6
$ws.Range("J88").Value = 1556.5714
$ws.Range("L88").Value = 1556.5714
$ws.Range("N88").Value = -2368.5714
$ws.Range("H91").Value = 1389.6
$ws.Range("J91").Value = 1556.5714
$ws.Range("L91").Value = 1556.5714
$ws.Range("N91").Value = -4364.5714

# ---- Sheet BSM: 8 cell updates ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5387
$ws.Range("J20").Value = 6330.5
$ws.Range("L20").Value = 6330.5
$ws.Range("N20").Value = -6824.5
$ws.Range("H130").Value = 61544
$ws.Range("J130").Value = 61544
$ws.Range("L130").Value = 61544
$ws.Range("N130").Value = -71584

# ---- Sheet CRP: 33 cell updates ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H58").Value = 3833.3333
$ws.Range("I58").Value = 2919.4
$ws.Range("J58").Value = 4486.143
$ws.Range("K58").Value = 2919.4
$ws.Range("L58").Value = 4486.143
$ws.Range("M58").Value = -2716.4
$ws.Range("N58").Value = -4892.143
$ws.Range("H132").Value = 2620.8518
$ws.Range("I132").Value = 2388.2104
$ws.Range("K132").Value = 7164.6312
$ws.Range("M132").Value = -4634.6312
$ws.Range("H134").Value = 2084.65
$ws.Range("I134").Value = 1674.3334
$ws.Range("J134").Value = 3315.6
$ws.Range("K134").Value = 5023.0002
$ws.Range("L134").Value = 9946.799999999999
$ws.Range("M134").Value = -2488.0002
$ws.Range("N134").Value = -15016.8
$ws.Range("H136").Value = 3833.3333
$ws.Range("I136").Value = 2919.4
$ws.Range("J136").Value = 4486.143
$ws.Range("K136").Value = 8758.200000000001
$ws.Range("L136").Value = 13458.429
$ws.Range("M136").Value = -6208.200000000001
$ws.Range("N136").Value = -18558.429

# ---- Sheet CUL: 35 cell updates ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 641.6
$ws.Range("I5").Value = 751
$ws.Range("K5").Value = 2253
$ws.Range("M5").Value = -2141
$ws.Range("H23").Value = 500099.5
$ws.Range("J23").Value = 500099.5
$ws.Range("L23").Value = 1500298.5
$ws.Range("N23").Value = -1500768.5
$ws.Range("H97").Value = 4965.6665
$ws.Range("J97").Value = 5813
$ws.Range("L97").Value = 17439
$ws.Range("N97").Value = -18431
$ws.Range("H121").Value = 792.6
$ws.Range("I121").Value = 382.33334
$ws.Range("J121").Value = 1408
$ws.Range("K121").Value = 1147.00002
$ws.Range("L121").Value = 4224
$ws.Range("M121").Value = 162.9999800000001
$ws.Range("N121").Value = -6844
$ws.Range("H135").Value = 641.6
$ws.Range("I135").Value = 751
$ws.Range("K135").Value = 6759
$ws.Range("M135").Value = -4224
$ws.Range("H137").Value = 4078.077
$ws.Range("I137").Value = 4021.6667
$ws.Range("K137").Value = 12065.0001
$ws.Range("M137").Value = -6965.000100000001
$ws.Range("H139").Value = 1904.5
$ws.Range("I139").Value = 1838.4445
$ws.Range("K139").Value = 5515.333500000001
$ws.Range("M139").Value = -375.3335000000006
$ws.Range("H140").Value = 4933.1304
$ws.Range("I140").Value = 2570.2222
$ws.Range("K140").Value = 7710.6666
$ws.Range("M140").Value = -2530.6666

# ---- Sheet GSM: 7 cell updates ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -19940

# ---- Sheet LTW: 48 cell updates ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8834
$ws.Range("I7").Value = 8834
$ws.Range("K7").Value = 8834
$ws.Range("M7").Value = -8722
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H61").Value = 15876801
$ws.Range("J61").Value = 4249
$ws.Range("L61").Value = 4249
$ws.Range("N61").Value = -4653
$ws.Range("H113").Value = 15876801
$ws.Range("J113").Value = 4249
$ws.Range("L113").Value = 4249
$ws.Range("N113").Value = -8589
$ws.Range("H122").Value = 4000
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 8834
$ws.Range("I126").Value = 8834
$ws.Range("K126").Value = 26502
$ws.Range("M126").Value = -24032
$ws.Range("H132").Value = 4122.625
$ws.Range("I132").Value = 3596.2
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10788.6
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -8258.599999999999
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 4569
$ws.Range("I136").Value = 4128.2666
$ws.Range("J136").Value = 6221.75
$ws.Range("K136").Value = 12384.7998
$ws.Range("L136").Value = 18665.25
$ws.Range("M136").Value = -9834.799800000001
$ws.Range("N136").Value = -23765.25

# ---- Sheet WVR: 15 cell updates ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 1534
$ws.Range("J101").Value = 1534
$ws.Range("L101").Value = 1534
$ws.Range("N101").Value = -8024
$ws.Range("H132").Value = 2338.1667
$ws.Range("I132").Value = 2338.9048
$ws.Range("J132").Value = 2333
$ws.Range("K132").Value = 7016.714399999999
$ws.Range("L132").Value = 6999
$ws.Range("M132").Value = -4486.714399999999
$ws.Range("N132").Value = -12059
$ws.Range("H136").Value = 2713.818
$ws.Range("I136").Value = 1626.2858
$ws.Range("K136").Value = 4878.857400000001
$ws.Range("M136").Value = -2328.857400000001

Write-Host "Applied 230 cell updates across 8 sheets"